{"js": "// Replace specific run texts identified by the diff.\n// Each edit: find the run whose text exactly equals `oldText` AND whose\n// enclosing paragraph starts with `oldText` (disambiguates against runs\n// where `oldText` is merely a substring of an already-updated, longer run).\nconst edits = [\n  { oldText: \"\u0432\u044a \u043f\u0440\u043e\u0441\u0442\u046b \u2192 \u0432\u044a + Acc. \u2192 \u0432\u044a & \u043f\u0440\u043e\u0441\u0442\u044a\", newText: \"\u0432\u044a \u043f\u0440\u043e\u0441\u0442\u046b \u2192 \u0432\u044a & \u043f\u0440\u043e\u0441\u0442\u044a\" },\n  { oldText: \"\u0442\u0432\u043e\u0440\ue205\u0442\ue205 & \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201\", newText: \"# \u0442\u0432\u043e\u0440\ue205\u0442\ue205 \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201 \u2192 \u0442\u0432\u043e\u0440\ue205\u0442\ue205 & \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201\" },\n  { oldText: \"\u0432\u044a \u043b\u0463\u043f\u043e\u0442\u046b \u2192 \u0432\u044a + Acc. \u2192 \u0432\u044a & \u043b\u0463\u043f\u043e\u0442\u0430\", newText: \"\u0432\u044a \u043b\u0463\u043f\u043e\u0442\u046b \u2192 \u0432\u044a & \u043b\u0463\u043f\u043e\u0442\u0430\" },\n  { oldText: \"\u043f\u043e \u043b\u0463\u043f\u043e\u0442\u0463 \u2192 \u043f\u043e + Dat. \u2192 \u043f\u043e & \u043b\u0463\u043f\u043e\u0442\u0430\", newText: \"\u043f\u043e \u043b\u0463\u043f\u043e\u0442\u0463 \u2192 \u043f\u043e & \u043b\u0463\u043f\u043e\u0442\u0430\" },\n  { oldText: \"\u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u044a & \ue20d\u0467\u0434\u044c\", newText: \"\u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u0430 \ue20d\u0467\u0434\u044c \u2192 \u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u044a & \ue20d\u0467\u0434\u044c\" },\n];\n\nconst body = context.document.body;\n\nfor (const edit of edits) {\n  const hits = body.search(edit.oldText, { matchCase: true });\n  hits.load('items');\n  await context.sync();\n\n  let applied = false;\n  for (const hit of hits.items) {\n    const para = hit.paragraphs.getFirst();\n    para.load('text');\n    await context.sync();\n    if (para.text.indexOf(edit.oldText) === 0) {\n      hit.insertText(edit.newText, Word.InsertLocation.replace);\n      applied = true;\n      break;\n    }\n  }\n  if (!applied) {\n    throw new Error('No unique match found for: ' + edit.oldText);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace specific run texts identified by the diff.\n# Each edit finds the paragraph whose text starts with `OldText` (this\n# disambiguates against paragraphs where `OldText` is merely a substring of\n# an already-updated, longer paragraph) and replaces just that leading run\n# text via a Find/Replace scoped to the paragraph's own Range.\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ OldText = \"\u0432\u044a \u043f\u0440\u043e\u0441\u0442\u046b \u2192 \u0432\u044a + Acc. \u2192 \u0432\u044a & \u043f\u0440\u043e\u0441\u0442\u044a\"; NewText = \"\u0432\u044a \u043f\u0440\u043e\u0441\u0442\u046b \u2192 \u0432\u044a & \u043f\u0440\u043e\u0441\u0442\u044a\" }\n    @{ OldText = \"\u0442\u0432\u043e\u0440\ue205\u0442\ue205 & \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201\"; NewText = \"# \u0442\u0432\u043e\u0440\ue205\u0442\ue205 \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201 \u2192 \u0442\u0432\u043e\u0440\ue205\u0442\ue205 & \ua641\u043d\u0430\u043c\u0435\u043d\ue205\ue201\" }\n    @{ OldText = \"\u0432\u044a \u043b\u0463\u043f\u043e\u0442\u046b \u2192 \u0432\u044a + Acc. \u2192 \u0432\u044a & \u043b\u0463\u043f\u043e\u0442\u0430\"; NewText = \"\u0432\u044a \u043b\u0463\u043f\u043e\u0442\u046b \u2192 \u0432\u044a & \u043b\u0463\u043f\u043e\u0442\u0430\" }\n    @{ OldText = \"\u043f\u043e \u043b\u0463\u043f\u043e\u0442\u0463 \u2192 \u043f\u043e + Dat. \u2192 \u043f\u043e & \u043b\u0463\u043f\u043e\u0442\u0430\"; NewText = \"\u043f\u043e \u043b\u0463\u043f\u043e\u0442\u0463 \u2192 \u043f\u043e & \u043b\u0463\u043f\u043e\u0442\u0430\" }\n    @{ OldText = \"\u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u044a & \ue20d\u0467\u0434\u044c\"; NewText = \"\u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u0430 \ue20d\u0467\u0434\u044c \u2192 \u0430\u0432\u0440\u0430\u0430\u043c\u043e\u0432\u044a & \ue20d\u0467\u0434\u044c\" }\n)\n\nforeach ($edit in $edits) {\n    $old = $edit.OldText\n    $new = $edit.NewText\n    $applied = $false\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($old)) {\n            $r = $p.Range\n            $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n            if (-not $found) {\n                throw \"Find.Execute failed for: $old\"\n            }\n            $applied = $true\n            break\n        }\n    }\n    if (-not $applied) {\n        throw \"No unique paragraph match found for: $old\"\n    }\n}\n"}
